$d = $word.ActiveDocument

$pairs = @(
    @{old="853÷7=121, 6"; new="823÷2=411, 1"},
    @{old="778÷8=97, 2"; new="859÷5=171, 4"},
    @{old="187÷4=46, 3"; new="617÷8=77, 1"},
    @{old="282÷6=47, 0"; new="420÷8=52, 4"},
    @{old="959÷8=119, 7"; new="627÷4=156, 3"},
    @{old="822÷4=205, 2"; new="707÷7=101, 0"},
    @{old="488÷2=244, 0"; new="619÷2=309, 1"},
    @{old="603÷3=201, 0"; new="960÷4=240, 0"},
    @{old="720÷4=180, 0"; new="115÷6=19, 1"},
    @{old="627÷3=209, 0"; new="970÷2=485, 0"},
    @{old="689÷4=172, 1"; new="628÷5=125, 3"},
    @{old="445÷8=55, 5"; new="284÷8=35, 4"},
    @{old="779÷5=155, 4"; new="758÷4=189, 2"},
    @{old="274÷2=137, 0"; new="359÷8=44, 7"},
    @{old="798÷3=266, 0"; new="454÷6=75, 4"},
    @{old="300÷5=60, 0"; new="113÷9=12, 5"},
    @{old="437÷6=72, 5"; new="692÷3=230, 2"},
    @{old="231÷8=28, 7"; new="543÷7=77, 4"},
    @{old="661÷6=110, 1"; new="692÷6=115, 2"},
    @{old="687÷6=114, 3"; new="605÷8=75, 5"},
    @{old="922÷7=131, 5"; new="576÷6=96, 0"},
    @{old="608÷2=304, 0"; new="284÷4=71, 0"},
    @{old="343÷3=114, 1"; new="411÷2=205, 1"},
    @{old="957÷5=191, 2"; new="919÷3=306, 1"},
    @{old="954÷6=159, 0"; new="261÷9=29, 0"}
)

foreach ($pair in $pairs) {
    $range = $d.Content
    $range.Find.Execute($pair.old, $true, $false, $false, $false, $false, $true, 1, $false, $pair.new, 2)
}
